$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.450"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06357"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'6.560"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8204"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.01412"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.1677"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08691"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03670"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03210"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "'3.709"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001644"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04751"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006232"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006275"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Value = "'3.784"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.270"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3354"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.04772"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007150"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.004506"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01167"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006807"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.004224"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"
$ws.Range("D50").Value = "'0.01242"
$ws.Range("D50").Style = "Normal"
